$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix sprint data: rows 2-11 (swap no_issue_added/no_issue_removed into no_issue_starttime/
# no_issue_todo pattern) and append rows 12-21 for the new sprint (d74cdee3-...).

# Row 2
$ws.Cells.Item(2, 1).Value = "253ac332-b8e2-43a7-bbc7-d673ddabc733"
$ws.Cells.Item(2, 2).Value = "bbfddaeb-dea5-4178-b067-e0fa8812c338"
$ws.Cells.Item(2, 3).Value = 13.0
$ws.Cells.Item(2, 4).Value = 10.0
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 0.0
$ws.Cells.Item(2, 7).Value = 0.0
$ws.Cells.Item(2, 8).Value = 0.0
$ws.Cells.Item(2, 9).Value = 3.0
$ws.Cells.Item(2, 10).Value = 0.0
$ws.Cells.Item(2, 11).Value = 3.0

# Row 3
$ws.Cells.Item(3, 1).Value = "3438c873-bc50-42c4-a598-1ed31761fe23"
$ws.Cells.Item(3, 2).Value = "bbfddaeb-dea5-4178-b067-e0fa8812c338"
$ws.Cells.Item(3, 3).Value = 13.0
$ws.Cells.Item(3, 4).Value = 10.0
$ws.Cells.Item(3, 5).Value = 4.0
$ws.Cells.Item(3, 6).Value = 0.0
$ws.Cells.Item(3, 7).Value = 0.0
$ws.Cells.Item(3, 8).Value = 0.0
$ws.Cells.Item(3, 9).Value = 4.0
$ws.Cells.Item(3, 10).Value = 0.0
$ws.Cells.Item(3, 11).Value = 3.0

# Row 4
$ws.Cells.Item(4, 1).Value = "80bc4a94-2ad4-44e8-bc20-8e76a44864e5"
$ws.Cells.Item(4, 2).Value = "bbfddaeb-dea5-4178-b067-e0fa8812c338"
$ws.Cells.Item(4, 3).Value = 13.0
$ws.Cells.Item(4, 4).Value = 10.0
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 0.0
$ws.Cells.Item(4, 7).Value = 0.0
$ws.Cells.Item(4, 8).Value = 0.0
$ws.Cells.Item(4, 9).Value = 3.0
$ws.Cells.Item(4, 10).Value = 0.0
$ws.Cells.Item(4, 11).Value = 1.0

# Row 5
$ws.Cells.Item(5, 1).Value = "89427d93-69a9-4187-baac-7a3105624b8b"
$ws.Cells.Item(5, 2).Value = "bbfddaeb-dea5-4178-b067-e0fa8812c338"
$ws.Cells.Item(5, 3).Value = 13.0
$ws.Cells.Item(5, 4).Value = 10.0
$ws.Cells.Item(5, 5).Value = 4.0
$ws.Cells.Item(5, 6).Value = 0.0
$ws.Cells.Item(5, 7).Value = 0.0
$ws.Cells.Item(5, 8).Value = 0.0
$ws.Cells.Item(5, 9).Value = 3.0
$ws.Cells.Item(5, 10).Value = 1.0
$ws.Cells.Item(5, 11).Value = 3.0

# Row 6
$ws.Cells.Item(6, 1).Value = "a5bc1288-fe8b-41c1-8993-152dbe7cb3fa"
$ws.Cells.Item(6, 2).Value = "bbfddaeb-dea5-4178-b067-e0fa8812c338"
$ws.Cells.Item(6, 3).Value = 13.0
$ws.Cells.Item(6, 4).Value = 10.0
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 0.0
$ws.Cells.Item(6, 7).Value = 0.0
$ws.Cells.Item(6, 8).Value = 0.0
$ws.Cells.Item(6, 9).Value = 2.0
$ws.Cells.Item(6, 10).Value = 1.0
$ws.Cells.Item(6, 11).Value = 1.0

# Row 7
$ws.Cells.Item(7, 1).Value = "a8eb1192-3e81-42e8-9f69-8483736af936"
$ws.Cells.Item(7, 2).Value = "bbfddaeb-dea5-4178-b067-e0fa8812c338"
$ws.Cells.Item(7, 3).Value = 13.0
$ws.Cells.Item(7, 4).Value = 10.0
$ws.Cells.Item(7, 5).Value = 4.0
$ws.Cells.Item(7, 6).Value = 0.0
$ws.Cells.Item(7, 7).Value = 0.0
$ws.Cells.Item(7, 8).Value = 0.0
$ws.Cells.Item(7, 9).Value = 4.0
$ws.Cells.Item(7, 10).Value = 0.0
$ws.Cells.Item(7, 11).Value = 1.0

# Row 8
$ws.Cells.Item(8, 1).Value = "cfcba21d-966a-43b8-a1a8-9751e05ed7c5"
$ws.Cells.Item(8, 2).Value = "bbfddaeb-dea5-4178-b067-e0fa8812c338"
$ws.Cells.Item(8, 3).Value = 13.0
$ws.Cells.Item(8, 4).Value = 10.0
$ws.Cells.Item(8, 5).Value = 4.0
$ws.Cells.Item(8, 6).Value = 0.0
$ws.Cells.Item(8, 7).Value = 0.0
$ws.Cells.Item(8, 8).Value = 0.0
$ws.Cells.Item(8, 9).Value = 4.0
$ws.Cells.Item(8, 10).Value = 0.0
$ws.Cells.Item(8, 11).Value = 3.0

# Row 9
$ws.Cells.Item(9, 1).Value = "d909ec23-e3b8-473c-8b53-b0dd10e5cb4a"
$ws.Cells.Item(9, 2).Value = "bbfddaeb-dea5-4178-b067-e0fa8812c338"
$ws.Cells.Item(9, 3).Value = 13.0
$ws.Cells.Item(9, 4).Value = 10.0
$ws.Cells.Item(9, 5).Value = 4.0
$ws.Cells.Item(9, 6).Value = 0.0
$ws.Cells.Item(9, 7).Value = 0.0
$ws.Cells.Item(9, 8).Value = 0.0
$ws.Cells.Item(9, 9).Value = 3.0
$ws.Cells.Item(9, 10).Value = 1.0
$ws.Cells.Item(9, 11).Value = 2.0

# Row 10
$ws.Cells.Item(10, 1).Value = "dc9af886-7862-4815-9876-037e5440de12"
$ws.Cells.Item(10, 2).Value = "bbfddaeb-dea5-4178-b067-e0fa8812c338"
$ws.Cells.Item(10, 3).Value = 13.0
$ws.Cells.Item(10, 4).Value = 10.0
$ws.Cells.Item(10, 5).Value = 2.0
$ws.Cells.Item(10, 6).Value = 0.0
$ws.Cells.Item(10, 7).Value = 0.0
$ws.Cells.Item(10, 8).Value = 0.0
$ws.Cells.Item(10, 9).Value = 2.0
$ws.Cells.Item(10, 10).Value = 0.0
$ws.Cells.Item(10, 11).Value = 1.0

# Row 11
$ws.Cells.Item(11, 1).Value = "f27d5a34-ecc1-4c32-8fd9-44615cbbed19"
$ws.Cells.Item(11, 2).Value = "bbfddaeb-dea5-4178-b067-e0fa8812c338"
$ws.Cells.Item(11, 3).Value = 13.0
$ws.Cells.Item(11, 4).Value = 10.0
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 0.0
$ws.Cells.Item(11, 7).Value = 0.0
$ws.Cells.Item(11, 8).Value = 0.0
$ws.Cells.Item(11, 9).Value = 3.0
$ws.Cells.Item(11, 10).Value = 0.0
$ws.Cells.Item(11, 11).Value = 3.0

# Row 12
$ws.Cells.Item(12, 1).Value = "253ac332-b8e2-43a7-bbc7-d673ddabc733"
$ws.Cells.Item(12, 2).Value = "d74cdee3-a4fc-44b4-bd90-1e53bec5cad1"
$ws.Cells.Item(12, 3).Value = 13.0
$ws.Cells.Item(12, 4).Value = 10.0
$ws.Cells.Item(12, 5).Value = 4.0
$ws.Cells.Item(12, 6).Value = 0.0
$ws.Cells.Item(12, 7).Value = 0.0
$ws.Cells.Item(12, 8).Value = 1.0
$ws.Cells.Item(12, 9).Value = 3.0
$ws.Cells.Item(12, 10).Value = 0.0
$ws.Cells.Item(12, 11).Value = 3.0

# Row 13
$ws.Cells.Item(13, 1).Value = "3438c873-bc50-42c4-a598-1ed31761fe23"
$ws.Cells.Item(13, 2).Value = "d74cdee3-a4fc-44b4-bd90-1e53bec5cad1"
$ws.Cells.Item(13, 3).Value = 13.0
$ws.Cells.Item(13, 4).Value = 10.0
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 0.0
$ws.Cells.Item(13, 7).Value = 0.0
$ws.Cells.Item(13, 8).Value = 0.0
$ws.Cells.Item(13, 9).Value = 3.0
$ws.Cells.Item(13, 10).Value = 0.0
$ws.Cells.Item(13, 11).Value = 3.0

# Row 14
$ws.Cells.Item(14, 1).Value = "80bc4a94-2ad4-44e8-bc20-8e76a44864e5"
$ws.Cells.Item(14, 2).Value = "d74cdee3-a4fc-44b4-bd90-1e53bec5cad1"
$ws.Cells.Item(14, 3).Value = 13.0
$ws.Cells.Item(14, 4).Value = 10.0
$ws.Cells.Item(14, 5).Value = 3.0
$ws.Cells.Item(14, 6).Value = 0.0
$ws.Cells.Item(14, 7).Value = 0.0
$ws.Cells.Item(14, 8).Value = 0.0
$ws.Cells.Item(14, 9).Value = 2.0
$ws.Cells.Item(14, 10).Value = 1.0
$ws.Cells.Item(14, 11).Value = 3.0

# Row 15
$ws.Cells.Item(15, 1).Value = "89427d93-69a9-4187-baac-7a3105624b8b"
$ws.Cells.Item(15, 2).Value = "d74cdee3-a4fc-44b4-bd90-1e53bec5cad1"
$ws.Cells.Item(15, 3).Value = 13.0
$ws.Cells.Item(15, 4).Value = 10.0
$ws.Cells.Item(15, 5).Value = 4.0
$ws.Cells.Item(15, 6).Value = 0.0
$ws.Cells.Item(15, 7).Value = 0.0
$ws.Cells.Item(15, 8).Value = 1.0
$ws.Cells.Item(15, 9).Value = 2.0
$ws.Cells.Item(15, 10).Value = 1.0
$ws.Cells.Item(15, 11).Value = 1.0

# Row 16
$ws.Cells.Item(16, 1).Value = "a5bc1288-fe8b-41c1-8993-152dbe7cb3fa"
$ws.Cells.Item(16, 2).Value = "d74cdee3-a4fc-44b4-bd90-1e53bec5cad1"
$ws.Cells.Item(16, 3).Value = 13.0
$ws.Cells.Item(16, 4).Value = 10.0
$ws.Cells.Item(16, 5).Value = 4.0
$ws.Cells.Item(16, 6).Value = 0.0
$ws.Cells.Item(16, 7).Value = 0.0
$ws.Cells.Item(16, 8).Value = 1.0
$ws.Cells.Item(16, 9).Value = 1.0
$ws.Cells.Item(16, 10).Value = 2.0
$ws.Cells.Item(16, 11).Value = 2.0

# Row 17
$ws.Cells.Item(17, 1).Value = "a8eb1192-3e81-42e8-9f69-8483736af936"
$ws.Cells.Item(17, 2).Value = "d74cdee3-a4fc-44b4-bd90-1e53bec5cad1"
$ws.Cells.Item(17, 3).Value = 13.0
$ws.Cells.Item(17, 4).Value = 10.0
$ws.Cells.Item(17, 5).Value = 3.0
$ws.Cells.Item(17, 6).Value = 0.0
$ws.Cells.Item(17, 7).Value = 0.0
$ws.Cells.Item(17, 8).Value = 0.0
$ws.Cells.Item(17, 9).Value = 1.0
$ws.Cells.Item(17, 10).Value = 2.0
$ws.Cells.Item(17, 11).Value = 2.0

# Row 18
$ws.Cells.Item(18, 1).Value = "cfcba21d-966a-43b8-a1a8-9751e05ed7c5"
$ws.Cells.Item(18, 2).Value = "d74cdee3-a4fc-44b4-bd90-1e53bec5cad1"
$ws.Cells.Item(18, 3).Value = 13.0
$ws.Cells.Item(18, 4).Value = 10.0
$ws.Cells.Item(18, 5).Value = 5.0
$ws.Cells.Item(18, 6).Value = 0.0
$ws.Cells.Item(18, 7).Value = 0.0
$ws.Cells.Item(18, 8).Value = 1.0
$ws.Cells.Item(18, 9).Value = 4.0
$ws.Cells.Item(18, 10).Value = 0.0
$ws.Cells.Item(18, 11).Value = 2.0

# Row 19
$ws.Cells.Item(19, 1).Value = "d909ec23-e3b8-473c-8b53-b0dd10e5cb4a"
$ws.Cells.Item(19, 2).Value = "d74cdee3-a4fc-44b4-bd90-1e53bec5cad1"
$ws.Cells.Item(19, 3).Value = 13.0
$ws.Cells.Item(19, 4).Value = 10.0
$ws.Cells.Item(19, 5).Value = 5.0
$ws.Cells.Item(19, 6).Value = 0.0
$ws.Cells.Item(19, 7).Value = 0.0
$ws.Cells.Item(19, 8).Value = 1.0
$ws.Cells.Item(19, 9).Value = 4.0
$ws.Cells.Item(19, 10).Value = 0.0
$ws.Cells.Item(19, 11).Value = 3.0

# Row 20
$ws.Cells.Item(20, 1).Value = "dc9af886-7862-4815-9876-037e5440de12"
$ws.Cells.Item(20, 2).Value = "d74cdee3-a4fc-44b4-bd90-1e53bec5cad1"
$ws.Cells.Item(20, 3).Value = 13.0
$ws.Cells.Item(20, 4).Value = 10.0
$ws.Cells.Item(20, 5).Value = 3.0
$ws.Cells.Item(20, 6).Value = 0.0
$ws.Cells.Item(20, 7).Value = 0.0
$ws.Cells.Item(20, 8).Value = 0.0
$ws.Cells.Item(20, 9).Value = 3.0
$ws.Cells.Item(20, 10).Value = 0.0
$ws.Cells.Item(20, 11).Value = 1.0

# Row 21
$ws.Cells.Item(21, 1).Value = "f27d5a34-ecc1-4c32-8fd9-44615cbbed19"
$ws.Cells.Item(21, 2).Value = "d74cdee3-a4fc-44b4-bd90-1e53bec5cad1"
$ws.Cells.Item(21, 3).Value = 13.0
$ws.Cells.Item(21, 4).Value = 10.0
$ws.Cells.Item(21, 5).Value = 3.0
$ws.Cells.Item(21, 6).Value = 0.0
$ws.Cells.Item(21, 7).Value = 0.0
$ws.Cells.Item(21, 8).Value = 0.0
$ws.Cells.Item(21, 9).Value = 3.0
$ws.Cells.Item(21, 10).Value = 0.0
$ws.Cells.Item(21, 11).Value = 2.0

